# Update "想去人数" (wishlist/interested count) figures that changed
# between the previous and newly generated data snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" -----------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 1212
$ws1.Range("F5").Value  = 36
$ws1.Range("F12").Value = 11717
$ws1.Range("F14").Value = 1386
$ws1.Range("F15").Value = 4654
$ws1.Range("F16").Value = 465
$ws1.Range("F19").Value = 70

# --- Sheet "全部类型" ---------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 1212
$ws4.Range("F6").Value  = 36
$ws4.Range("F15").Value = 11717
$ws4.Range("F17").Value = 1386
$ws4.Range("F18").Value = 4654
$ws4.Range("F19").Value = 465
$ws4.Range("F22").Value = 70
